# Append the next day's GSC export row to the "Chart" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Find the last used row in column A and append right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Write the date as literal text (not an auto-converted date serial) by
# routing it through a text-forcing formula, then freezing the result to a
# plain value via copy/paste-special so no formula or stray number format
# is left behind on the cell.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Formula = "=T(""2025-12-28"")"
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 28
